{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of: async (context) => { ... }\n\nconst replacements = [\n  [\"2024-02-04 Sunday\", \"2024-02-05 Monday\"],\n  [\"124\\u00F74=31, 0\", \"260\\u00F76=43, 2\"],\n  [\"740\\u00F78=92, 4\", \"120\\u00F75=24, 0\"],\n  [\"876\\u00F75=175, 1\", \"433\\u00F72=216, 1\"],\n  [\"501\\u00F75=100, 1\", \"462\\u00F73=154, 0\"],\n  [\"406\\u00F74=101, 2\", \"202\\u00F73=67, 1\"],\n  [\"764\\u00F75=152, 4\", \"878\\u00F74=219, 2\"],\n  [\"259\\u00F77=37, 0\", \"353\\u00F79=39, 2\"],\n  [\"525\\u00F77=75, 0\", \"286\\u00F72=143, 0\"],\n  [\"202\\u00F76=33, 4\", \"320\\u00F72=160, 0\"],\n  [\"754\\u00F73=251, 1\", \"179\\u00F75=35, 4\"],\n  [\"420\\u00F76=70, 0\", \"197\\u00F78=24, 5\"],\n  [\"151\\u00F73=50, 1\", \"235\\u00F73=78, 1\"],\n  [\"960\\u00F72=480, 0\", \"783\\u00F78=97, 7\"],\n  [\"565\\u00F77=80, 5\", \"372\\u00F78=46, 4\"],\n  [\"188\\u00F79=20, 8\", \"994\\u00F77=142, 0\"],\n  [\"532\\u00F79=59, 1\", \"374\\u00F76=62, 2\"],\n  [\"208\\u00F74=52, 0\", \"485\\u00F75=97, 0\"],\n  [\"402\\u00F78=50, 2\", \"947\\u00F76=157, 5\"],\n  [\"908\\u00F72=454, 0\", \"869\\u00F78=108, 5\"],\n  [\"965\\u00F72=482, 1\", \"613\\u00F75=122, 3\"],\n  [\"550\\u00F73=183, 1\", \"705\\u00F77=100, 5\"],\n  [\"946\\u00F77=135, 1\", \"651\\u00F77=93, 0\"],\n  [\"118\\u00F79=13, 1\", \"475\\u00F72=237, 1\"],\n  [\"637\\u00F78=79, 5\", \"162\\u00F77=23, 1\"],\n  [\"738\\u00F77=105, 3\", \"402\\u00F78=50, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace only the first occurrence (each source string is unique in\n  // this document), mirroring the diff which edits exactly one run per text.\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $app / $doc resolve to the running application / ActiveDocument.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-04 Sunday\", \"2024-02-05 Monday\"),\n    @(\"124\u00f74=31, 0\", \"260\u00f76=43, 2\"),\n    @(\"740\u00f78=92, 4\", \"120\u00f75=24, 0\"),\n    @(\"876\u00f75=175, 1\", \"433\u00f72=216, 1\"),\n    @(\"501\u00f75=100, 1\", \"462\u00f73=154, 0\"),\n    @(\"406\u00f74=101, 2\", \"202\u00f73=67, 1\"),\n    @(\"764\u00f75=152, 4\", \"878\u00f74=219, 2\"),\n    @(\"259\u00f77=37, 0\", \"353\u00f79=39, 2\"),\n    @(\"525\u00f77=75, 0\", \"286\u00f72=143, 0\"),\n    @(\"202\u00f76=33, 4\", \"320\u00f72=160, 0\"),\n    @(\"754\u00f73=251, 1\", \"179\u00f75=35, 4\"),\n    @(\"420\u00f76=70, 0\", \"197\u00f78=24, 5\"),\n    @(\"151\u00f73=50, 1\", \"235\u00f73=78, 1\"),\n    @(\"960\u00f72=480, 0\", \"783\u00f78=97, 7\"),\n    @(\"565\u00f77=80, 5\", \"372\u00f78=46, 4\"),\n    @(\"188\u00f79=20, 8\", \"994\u00f77=142, 0\"),\n    @(\"532\u00f79=59, 1\", \"374\u00f76=62, 2\"),\n    @(\"208\u00f74=52, 0\", \"485\u00f75=97, 0\"),\n    @(\"402\u00f78=50, 2\", \"947\u00f76=157, 5\"),\n    @(\"908\u00f72=454, 0\", \"869\u00f78=108, 5\"),\n    @(\"965\u00f72=482, 1\", \"613\u00f75=122, 3\"),\n    @(\"550\u00f73=183, 1\", \"705\u00f77=100, 5\"),\n    @(\"946\u00f77=135, 1\", \"651\u00f77=93, 0\"),\n    @(\"118\u00f79=13, 1\", \"475\u00f72=237, 1\"),\n    @(\"637\u00f78=79, 5\", \"162\u00f77=23, 1\"),\n    @(\"738\u00f77=105, 3\", \"402\u00f78=50, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Signature: FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    # MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    # ReplaceWith, Replace (2 = wdReplaceAll)\n    $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        0,\n        $false,\n        $newText,\n        2\n    )\n}\n\n$d.Save()\n"}
